$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1903669724770642
$ws.Range("C2").Value = 0.5619266055045872
$ws.Range("J2").Value = 0.006880733944954129
$ws.Range("P2").Value = 0.1605504587155963
$ws.Range("S2").Value = 0.08027522935779817
$ws.Range("B3").Value = 0.02777777777777778
$ws.Range("C3").Value = 0.01984126984126984
$ws.Range("J3").Value = 0.007936507936507936
$ws.Range("P3").Value = 0.753968253968254
$ws.Range("S3").Value = 0.1904761904761905
$ws.Range("J4").Value = 0.04411764705882353
$ws.Range("O4").Value = 0.01470588235294118
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2647058823529412
$ws.Range("B6").Value = 0.09278350515463918
$ws.Range("D6").Value = 0.01030927835051546
$ws.Range("F6").Value = 0.04639175257731959
$ws.Range("J6").Value = 0.2474226804123711
$ws.Range("O6").Value = 0.04639175257731959
$ws.Range("Q6").Value = 0.154639175257732
$ws.Range("R6").Value = 0.07216494845360824
$ws.Range("S6").Value = 0.3298969072164948
$ws.Range("B7").Value = 0.1784037558685446
$ws.Range("D7").Value = 0.07042253521126761
$ws.Range("F7").Value = 0.02816901408450704
$ws.Range("J7").Value = 0.09859154929577464
$ws.Range("O7").Value = 0.009389671361502348
$ws.Range("Q7").Value = 0.1643192488262911
$ws.Range("R7").Value = 0.09389671361502347
$ws.Range("S7").Value = 0.3568075117370892
$ws.Range("B8").Value = 0.1164241164241164
$ws.Range("D8").Value = 0.02910602910602911
$ws.Range("F8").Value = 0.04158004158004158
$ws.Range("J8").Value = 0.07692307692307693
$ws.Range("O8").Value = 0.0103950103950104
$ws.Range("Q8").Value = 0.1891891891891892
$ws.Range("R8").Value = 0.08731808731808732
$ws.Range("S8").Value = 0.4490644490644491
$ws.Range("B9").Value = 0.1187739463601533
$ws.Range("D9").Value = 0.02298850574712644
$ws.Range("E9").Value = 0.003831417624521073
$ws.Range("F9").Value = 0.07662835249042145
$ws.Range("J9").Value = 0.06513409961685823
$ws.Range("O9").Value = 0.003831417624521073
$ws.Range("Q9").Value = 0.1762452107279693
$ws.Range("R9").Value = 0.103448275862069
$ws.Range("S9").Value = 0.4291187739463602
$ws.Range("B10").Value = 0.1413373860182371
$ws.Range("D10").Value = 0.02583586626139818
$ws.Range("E10").Value = 0.001519756838905775
$ws.Range("F10").Value = 0.05395136778115502
$ws.Range("J10").Value = 0.07674772036474165
$ws.Range("O10").Value = 0.01367781155015198
$ws.Range("Q10").Value = 0.2264437689969605
$ws.Range("R10").Value = 0.0851063829787234
$ws.Range("S10").Value = 0.3753799392097265
$ws.Range("G11").Value = 0.1577380952380952
$ws.Range("J11").Value = 0.07738095238095238
$ws.Range("K11").Value = 0.2172619047619048
$ws.Range("L11").Value = 0.5357142857142857
$ws.Range("S11").Value = 0.0119047619047619
$ws.Range("G12").Value = 0.7248677248677249
$ws.Range("J12").Value = 0.1904761904761905
$ws.Range("K12").Value = 0.02645502645502645
$ws.Range("L12").Value = 0.04232804232804233
$ws.Range("S12").Value = 0.01587301587301587
$ws.Range("G13").Value = 0.6170212765957447
$ws.Range("J13").Value = 0.3404255319148936
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.01388888888888889
$ws.Range("H15").Value = 0.1759259259259259
$ws.Range("I15").Value = 0.1064814814814815
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.04629629629629629
$ws.Range("M15").Value = 0.009259259259259259
$ws.Range("O15").Value = 0.02314814814814815
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.01689189189189189
$ws.Range("H16").Value = 0.1587837837837838
$ws.Range("I16").Value = 0.06418918918918919
$ws.Range("J16").Value = 0.4324324324324325
$ws.Range("K16").Value = 0.1081081081081081
$ws.Range("M16").Value = 0.01013513513513514
$ws.Range("N16").Value = 0.003378378378378379
$ws.Range("O16").Value = 0.04391891891891892
$ws.Range("S16").Value = 0.1621621621621622
$ws.Range("F17").Value = 0.02235772357723577
$ws.Range("H17").Value = 0.1585365853658537
$ws.Range("I17").Value = 0.1158536585365854
$ws.Range("J17").Value = 0.443089430894309
$ws.Range("K17").Value = 0.09552845528455285
$ws.Range("M17").Value = 0.01016260162601626
$ws.Range("O17").Value = 0.04471544715447155
$ws.Range("S17").Value = 0.1097560975609756
$ws.Range("F18").Value = 0.04225352112676056
$ws.Range("H18").Value = 0.1502347417840376
$ws.Range("I18").Value = 0.1126760563380282
$ws.Range("J18").Value = 0.4178403755868544
$ws.Range("K18").Value = 0.0892018779342723
$ws.Range("M18").Value = 0.009389671361502348
$ws.Range("N18").Value = 0.004694835680751174
$ws.Range("O18").Value = 0.06103286384976526
$ws.Range("S18").Value = 0.1126760563380282
$ws.Range("F19").Value = 0.01049685094471658
$ws.Range("H19").Value = 0.2008397480755773
$ws.Range("I19").Value = 0.09657102869139259
$ws.Range("J19").Value = 0.3519944016794961
$ws.Range("K19").Value = 0.1042687193841847
$ws.Range("M19").Value = 0.02659202239328201
$ws.Range("N19").Value = 0.0006997900629811056
$ws.Range("O19").Value = 0.06927921623512946
$ws.Range("S19").Value = 0.13925822253324
